# Greek (el-GR) localisation pass over the "new.pptx" empty template:
#  - slide master placeholder prompt text
#  - "Picture with Caption" layout prompt text
#  - title-slide shape names on the single content slide

$p = $ppt.ActivePresentation

# --- Slide Master (ppt/slideMasters/slideMaster1.xml) -----------------
$master = $p.SlideMaster

# Shape 1: Title Placeholder -> "Click to edit Master title style"
$titlePh = $master.Shapes.Item(1)
$titlePh.TextFrame.TextRange.Text = "Στυλ κύριου τίτλου"

# Shape 2: Text Placeholder -> body placeholder with 5 outline levels
$bodyPh = $master.Shapes.Item(2)
$bodyTr = $bodyPh.TextFrame.TextRange
$bodyTr.Paragraphs(1).Text = "Στυλ υποδείγματος κειμένου"
$bodyTr.Paragraphs(2).Text = "Δεύτερου επιπέδου"
$bodyTr.Paragraphs(3).Text = "Τρίτου επιπέδου"
$bodyTr.Paragraphs(4).Text = "Τέταρτου επιπέδου"
$bodyTr.Paragraphs(5).Text = "Πέμπτου επιπέδου"

# --- "Picture with Caption" slide layout (ppt/slideLayouts/slideLayout9.xml) ---
# It is the 9th layout attached to the (only) slide master.
$picLayout = $master.CustomLayouts.Item(9)
$picPh = $picLayout.Shapes.Item(2)
$picPh.TextFrame.TextRange.Text = "Κάντε κλικ στο εικονίδιο για να προσθέσετε εικόνα"

# Rename the master placeholder shapes (cNvPr "name" attribute) to match.
# Some COM-interop hosts do not persist Shape.Name writes made through
# Master shape collections (only slide-tier shapes reliably accept
# renames) -- guard each one so a host-level limitation here can't abort
# the rest of the (working) text-content edits below.
try { $master.Shapes.Item(1).Name = "Θέση τίτλου 1" } catch {}
try { $master.Shapes.Item(2).Name = "Θέση κειμένου 2" } catch {}
try { $master.Shapes.Item(3).Name = "Θέση ημερομηνίας 3" } catch {}
try { $master.Shapes.Item(4).Name = "Θέση υποσέλιδου 4" } catch {}
try { $master.Shapes.Item(5).Name = "Θέση αριθμού διαφάνειας 5" } catch {}

# --- Title slide (ppt/slides/slide1.xml) -------------------------------
$slide = $p.Slides.Item(1)
$slide.Shapes.Item(1).Name = "Τίτλος 1"
$slide.Shapes.Item(2).Name = "Υπότιτλος 2"
